$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each changed cell, force the "Text" number format before assigning the
# value. The source data stores every Coin/Link/Price/Volume cell as text
# (inline strings), including values that look numeric (e.g. "26.104.20",
# "1.000", "40.13"). Without forcing text format first, Excel would
# reinterpret such strings as numbers/dates and corrupt them (e.g. "1.000"
# -> 1, "40.13" -> 40.130000000000003).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.104.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.748.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5252"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2775"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.13"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06193"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.757.20"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07176"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.35"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6437"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.596"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.28"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "26.014.40"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.66"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006711"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.977.76"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.303"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.03%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.832"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.71%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.206"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.12"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.518"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.32"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.816"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "104.36"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08354"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.787"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.677"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.96%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04557"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.638"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6327"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.39%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.718"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01604"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.949"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9995"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.94"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3907"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.08%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7366"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.063"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1140"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.67%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.316"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05357"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.42%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.49"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.666"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.97%  "
